# Generate Report for Handoff
# Adds a new localization entry (bb2eb55e-b253-4e69-89ef-5182e71f4d83) to the
# Overview / zh-cn / de-de report sheets, pushing the ".localization-config"
# bookkeeping row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 9 gets what used to be in row 8 (the ".localization-config" bookkeeping row)
$ws1.Range("A9").Value = ".localization-config"
$ws1.Range("B9").Value = "Not to be localized"
$ws1.Range("C9").Value = "Not to be localized"

# Row 8 becomes the new handoff entry
$ws1.Range("A8").Value = "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md"
$ws1.Range("B8").Value = "Ready for handoff"
$ws1.Range("C8").Value = "Ready for handoff"

# Rebuild hyperlinks in final order (delete old collection to avoid stale/duplicate entries)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1116489707120925fb84edc7ef1a27b73277dadc/e2e/13241f23-34bc-4eff-b09c-39b84f281564.md", "", "", "13241f23-34bc-4eff-b09c-39b84f281564.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/6bac3023-21de-433a-b566-69529a3c67c3.md", "", "", "6bac3023-21de-433a-b566-69529a3c67c3.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md", "", "", "ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ffd23fa5bec35ff4f064ece4cec2788e52db5e1/e2e/ca734e4d-eca7-4b4b-bf24-686cbce93e15.md", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/e2e/058f6a44-efc6-4f84-98d4-8c23c5890d06.md", "", "", "058f6a44-efc6-4f84-98d4-8c23c5890d06.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1c40debd2cc2f34fd579cf1d29f2f81600806c0c/e2e/49ed7b88-ffcd-4894-879d-8a6b41754fa8.md", "", "", "49ed7b88-ffcd-4894-879d-8a6b41754fa8.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/a711cdd85ce8f26a2998e190cc888f39db1658ce/e2e/bb2eb55e-b253-4e69-89ef-5182e71f4d83.md", "", "", "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 9 gets what used to be in row 8
$ws2.Range("A9").Value = ".localization-config"
$ws2.Range("B9").Value = "Not to be localized"
$ws2.Range("D9").Value = "0001-01-01 00:00:00"
$ws2.Range("G9").Value = "0001-01-01 00:00:00"
$ws2.Range("H9").Value = "Ignored"

# Row 8 becomes the new handoff entry
$ws2.Range("A8").Value = "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md"
$ws2.Range("B8").Value = "Ready for handoff"
$ws2.Range("C8").Value = "bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.zh-cn.xlf"
$ws2.Range("D8").Value = "2016-03-01 09:03:17"
$ws2.Range("G8").Value = "0001-01-01 00:00:00"
$ws2.Range("H8").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1116489707120925fb84edc7ef1a27b73277dadc/e2e/13241f23-34bc-4eff-b09c-39b84f281564.md", "", "", "13241f23-34bc-4eff-b09c-39b84f281564.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6512858dfdba75ce9c0cdb52298461ef965d4033/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/13241f23-34bc-4eff-b09c-39b84f281564.48f2547baace9d553ea316ebff2f9ad538ae4e13.zh-cn.xlf", "", "", "13241f23-34bc-4eff-b09c-39b84f281564.48f2547baace9d553ea316ebff2f9ad538ae4e13.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/6bac3023-21de-433a-b566-69529a3c67c3.md", "", "", "6bac3023-21de-433a-b566-69529a3c67c3.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ac960e70070577ddbca3a33924bede3a13ebb73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6bac3023-21de-433a-b566-69529a3c67c3.fc448468024d66de0498c2cd582087e45e09ccbb.zh-cn.xlf", "", "", "6bac3023-21de-433a-b566-69529a3c67c3.fc448468024d66de0498c2cd582087e45e09ccbb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md", "", "", "ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ac960e70070577ddbca3a33924bede3a13ebb73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.7d5a12fd3af6eff7c675754d1beefc69d97d0613.zh-cn.xlf", "", "", "ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.7d5a12fd3af6eff7c675754d1beefc69d97d0613.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ffd23fa5bec35ff4f064ece4cec2788e52db5e1/e2e/ca734e4d-eca7-4b4b-bf24-686cbce93e15.md", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97c606f7e56372b45e5f9623643d885d770298c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.zh-cn.xlf", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/cbb159fc9d3809e760f45601efbc577c395f35e3/e2e/ca734e4d-eca7-4b4b-bf24-686cbce93e15.md", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/73465676ce28fa57fb07300c41ffebb0e58f1c2d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.zh-cn.xlf", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/e2e/058f6a44-efc6-4f84-98d4-8c23c5890d06.md", "", "", "058f6a44-efc6-4f84-98d4-8c23c5890d06.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df204dfd8e0ccb2de0a982ce33978c03d183a8b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/058f6a44-efc6-4f84-98d4-8c23c5890d06.bae27c346aff3dfcd9bc23cff8a1af51de882104.zh-cn.xlf", "", "", "058f6a44-efc6-4f84-98d4-8c23c5890d06.bae27c346aff3dfcd9bc23cff8a1af51de882104.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1c40debd2cc2f34fd579cf1d29f2f81600806c0c/e2e/49ed7b88-ffcd-4894-879d-8a6b41754fa8.md", "", "", "49ed7b88-ffcd-4894-879d-8a6b41754fa8.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ab0e65ddd6897d97282a6e6d223763af623adbf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/49ed7b88-ffcd-4894-879d-8a6b41754fa8.fa67e7d41d0de99e5250aa6fd0a1540c3bea163f.zh-cn.xlf", "", "", "49ed7b88-ffcd-4894-879d-8a6b41754fa8.fa67e7d41d0de99e5250aa6fd0a1540c3bea163f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/a711cdd85ce8f26a2998e190cc888f39db1658ce/e2e/bb2eb55e-b253-4e69-89ef-5182e71f4d83.md", "", "", "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a711cdd85ce8f26a2998e190cc888f39db1658ce/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.zh-cn.xlf", "", "", "bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 9 gets what used to be in row 8
$ws3.Range("A9").Value = ".localization-config"
$ws3.Range("B9").Value = "Not to be localized"
$ws3.Range("D9").Value = "0001-01-01 00:00:00"
$ws3.Range("G9").Value = "0001-01-01 00:00:00"
$ws3.Range("H9").Value = "Ignored"

# Row 8 becomes the new handoff entry
$ws3.Range("A8").Value = "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md"
$ws3.Range("B8").Value = "Ready for handoff"
$ws3.Range("C8").Value = "bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.de-de.xlf"
$ws3.Range("D8").Value = "2016-03-01 09:03:28"
$ws3.Range("G8").Value = "0001-01-01 00:00:00"
$ws3.Range("H8").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1116489707120925fb84edc7ef1a27b73277dadc/e2e/13241f23-34bc-4eff-b09c-39b84f281564.md", "", "", "13241f23-34bc-4eff-b09c-39b84f281564.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7256d713523ea00585243bf6c60f6924eb684cc2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/13241f23-34bc-4eff-b09c-39b84f281564.48f2547baace9d553ea316ebff2f9ad538ae4e13.de-de.xlf", "", "", "13241f23-34bc-4eff-b09c-39b84f281564.48f2547baace9d553ea316ebff2f9ad538ae4e13.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/6bac3023-21de-433a-b566-69529a3c67c3.md", "", "", "6bac3023-21de-433a-b566-69529a3c67c3.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71c1db93dafd155296eac524a898e8ca2f7dcd21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6bac3023-21de-433a-b566-69529a3c67c3.fc448468024d66de0498c2cd582087e45e09ccbb.de-de.xlf", "", "", "6bac3023-21de-433a-b566-69529a3c67c3.fc448468024d66de0498c2cd582087e45e09ccbb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d7e9b70457f1dc2e4b108f2976b321e2ba855e/e2e/ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md", "", "", "ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71c1db93dafd155296eac524a898e8ca2f7dcd21/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.7d5a12fd3af6eff7c675754d1beefc69d97d0613.de-de.xlf", "", "", "ac079ddf-571a-4acd-84fc-9ff1e2a6b9d2.7d5a12fd3af6eff7c675754d1beefc69d97d0613.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/5ffd23fa5bec35ff4f064ece4cec2788e52db5e1/e2e/ca734e4d-eca7-4b4b-bf24-686cbce93e15.md", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/767bf2ed2dfe6ae522a2cb3977531431f8a72a92/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.de-de.xlf", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/27ba2bd98389f4d5b7ff10c0c8b775c746154a41/e2e/ca734e4d-eca7-4b4b-bf24-686cbce93e15.md", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fce5919ac1297b47bed10dbe2bd03e89986ddec9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.de-de.xlf", "", "", "ca734e4d-eca7-4b4b-bf24-686cbce93e15.62bc92e7c4c992f8d388fc2723043dd8281f8190.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/e2e/058f6a44-efc6-4f84-98d4-8c23c5890d06.md", "", "", "058f6a44-efc6-4f84-98d4-8c23c5890d06.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b48b4db88b5b04e492b36cf59d67642e22a44db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/058f6a44-efc6-4f84-98d4-8c23c5890d06.bae27c346aff3dfcd9bc23cff8a1af51de882104.de-de.xlf", "", "", "058f6a44-efc6-4f84-98d4-8c23c5890d06.bae27c346aff3dfcd9bc23cff8a1af51de882104.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/1c40debd2cc2f34fd579cf1d29f2f81600806c0c/e2e/49ed7b88-ffcd-4894-879d-8a6b41754fa8.md", "", "", "49ed7b88-ffcd-4894-879d-8a6b41754fa8.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86ad613beaf36fb2e7fa77ce80b87571015fd669/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/49ed7b88-ffcd-4894-879d-8a6b41754fa8.fa67e7d41d0de99e5250aa6fd0a1540c3bea163f.de-de.xlf", "", "", "49ed7b88-ffcd-4894-879d-8a6b41754fa8.fa67e7d41d0de99e5250aa6fd0a1540c3bea163f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/a711cdd85ce8f26a2998e190cc888f39db1658ce/e2e/bb2eb55e-b253-4e69-89ef-5182e71f4d83.md", "", "", "bb2eb55e-b253-4e69-89ef-5182e71f4d83.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a711cdd85ce8f26a2998e190cc888f39db1658ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.de-de.xlf", "", "", "bb2eb55e-b253-4e69-89ef-5182e71f4d83.a711cdd85ce8f26a2998e190cc888f39db1658ce.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/a7b16235b044e0b981b2903bc9a5dabfc9736b42/.localization-config", "", "", ".localization-config") | Out-Null
